# Remove the trailing "Ver no Jupiter ..." / "© 2020 ..." footer block
# (and the blank paragraph that precedes it) that followed the
# "LOQ4095: Química Geral Experimental (Requisito fraco)" requirement line.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1

for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*LOQ4095*") {
        # first paragraph to remove is the one right after the LOQ4095 line
        $startIdx = $i + 1
    }
    if ($t -like "*Creative Commons Attribution*") {
        # last paragraph to remove is the copyright/footer line itself
        $endIdx = $i
    }
}

if ($startIdx -gt 0 -and $endIdx -ge $startIdx) {
    $startRange = $d.Paragraphs.Item($startIdx).Range
    $endRange = $d.Paragraphs.Item($endIdx).Range
    $rng = $d.Range($startRange.Start, $endRange.End)
    $rng.Delete()
}
